# setting up python example
#
# Insert a new "Graphical Programming Tools" slide right after
# "Common KLC Tasks" (i.e. at position 3), reusing the footer placeholder
# from the " Command Line" slide, then push the " Command Line" slide to
# the end of the deck.

$p = $ppt.ActivePresentation

# Grab the footer placeholder from the existing " Command Line" slide
# (currently slide 3) so the new slide can reuse the same footer
# placeholder/text.
$footerShape = $p.Slides.Item(3).Shapes.Item(2)
$footerShape.Copy()

# Add the new slide at index 3 using the same "Title, Content" layout used
# by the other content slides in this deck (layout 2 -> slideLayout2.xml).
$newSlide = $p.Slides.Add(3, 2)
$null = $newSlide.Shapes.Paste()

# --- Title -----------------------------------------------------------
$newSlide.Shapes.Item(1).Name = "Title 1"
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Graphical Programming Tools"

# --- Body bullets ------------------------------------------------------
$newSlide.Shapes.Item(2).Name = "Text Placeholder 2"
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Python"
$null = $body.InsertAfter("`rVSCode: ")
$null = $body.InsertAfter("`rJupyter notebooks: ")
$null = $body.InsertAfter("`rR")
$null = $body.InsertAfter("`rRStudio")
$null = $body.InsertAfter("`rJupyter notebooks")
$null = $body.InsertAfter("`rStata")
$null = $body.InsertAfter("`rXStata")

# Sub-bullets (editors/tools) sit one indent level deeper than the
# language headings.
$body.Paragraphs(2).IndentLevel = 2
$body.Paragraphs(3).IndentLevel = 2
$body.Paragraphs(5).IndentLevel = 2
$body.Paragraphs(6).IndentLevel = 2
$body.Paragraphs(8).IndentLevel = 2

# --- Footer --------------------------------------------------------
$newSlide.Shapes.Item(3).TextFrame.TextRange.Text = "Reproducibility Principles"

# Move the " Command Line" slide (now at index 4, after the insert) to the
# end of the deck.
$cmdLineSlide = $p.Slides.Item(4)
$cmdLineSlide.MoveTo($p.Slides.Count)
